# Re-append the 17 match rows for Shikhar Dhawan (Delhi Capitals) a second
# time at the bottom of the sheet (rows 19-35), mirroring the "json file
# code updated" re-export that duplicated the dataset. Column order is
# venue, date, result, ownTeam, oppTeam, batsman, totalRuns, totalBalls,
# total4s, total6s, sr - all stored as text, exactly like rows 2-18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(" Dubai (DSC)", " October 14 2020", "Capitals won by 13 runs", "Delhi Capitals", "Rajasthan Royals", "Shikhar Dhawan ", "57", "33", "6", "2", "172.72"),
    @(" Abu Dhabi", " October 11 2020", "Mumbai won by 5 wickets (with 2 balls remaining)", "Delhi Capitals", "Mumbai Indians", "Shikhar Dhawan ", "69", "52", "6", "1", "132.69"),
    @(" Abu Dhabi", " November 02 2020", "Capitals won by 6 wickets (with 6 balls remaining)", "Delhi Capitals", "Royal Challengers Bangalore", "Shikhar Dhawan ", "54", "41", "6", "0", "131.70"),
    @(" Dubai (DSC)", " November 05 2020", "Mumbai won by 57 runs", "Delhi Capitals", "Mumbai Indians", "Shikhar Dhawan ", "0", "2", "0", "0", "0.00"),
    @(" Dubai (DSC)", " October 27 2020", "Sunrisers won by 88 runs", "Delhi Capitals", "Sunrisers Hyderabad", "Shikhar Dhawan ", "0", "1", "0", "0", "0.00"),
    @(" Abu Dhabi", " September 29 2020", "Sunrisers won by 15 runs", "Delhi Capitals", "Sunrisers Hyderabad", "Shikhar Dhawan ", "34", "31", "4", "0", "109.67"),
    @(" Abu Dhabi", " October 24 2020", "KKR won by 59 runs", "Delhi Capitals", "Kolkata Knight Riders", "Shikhar Dhawan ", "6", "6", "1", "0", "100.00"),
    @(" Sharjah", " October 17 2020", "Capitals won by 5 wickets (with 1 ball remaining)", "Delhi Capitals", "Chennai Super Kings", "Shikhar Dhawan ", "101", "58", "14", "1", "174.13"),
    @(" Dubai (DSC)", " October 31 2020", "Mumbai won by 9 wickets (with 34 balls remaining)", "Delhi Capitals", "Mumbai Indians", "Shikhar Dhawan ", "0", "2", "0", "0", "0.00"),
    @(" Dubai (DSC)", " November 10 2020", "Mumbai won by 5 wickets (with 8 balls remaining)", "Delhi Capitals", "Mumbai Indians", "Shikhar Dhawan ", "15", "13", "3", "0", "115.38"),
    @(" Dubai (DSC)", " October 20 2020", "Kings XI won by 5 wickets (with 6 balls remaining)", "Delhi Capitals", "Kings XI Punjab", "Shikhar Dhawan ", "106", "61", "12", "3", "173.77"),
    @(" Dubai (DSC)", " October 05 2020", "Capitals won by 59 runs", "Delhi Capitals", "Royal Challengers Bangalore", "Shikhar Dhawan ", "32", "28", "3", "0", "114.28"),
    @(" Dubai (DSC)", " September 20 2020", "Match tied (Capitals won the one-over eliminator)", "Delhi Capitals", "Kings XI Punjab", "Shikhar Dhawan ", "0", "2", "0", "0", "0.00"),
    @(" Abu Dhabi", " November 08 2020", "Capitals won by 17 runs", "Delhi Capitals", "Sunrisers Hyderabad", "Shikhar Dhawan ", "78", "50", "6", "2", "156.00"),
    @(" Sharjah", " October 03 2020", "Capitals won by 18 runs", "Delhi Capitals", "Kolkata Knight Riders", "Shikhar Dhawan ", "26", "16", "2", "2", "162.50"),
    @(" Dubai (DSC)", " September 25 2020", "Capitals won by 44 runs", "Delhi Capitals", "Chennai Super Kings", "Shikhar Dhawan ", "35", "27", "3", "1", "129.62"),
    @(" Sharjah", " October 09 2020", "Capitals won by 46 runs", "Delhi Capitals", "Rajasthan Royals", "Shikhar Dhawan ", "5", "4", "1", "0", "125.00")
)

$startRow = 19
$lastRow = $startRow + $rows.Count - 1

# Columns G:K (totalRuns, totalBalls, total4s, total6s, sr) hold numeric
# looking values but must stay text, same as the existing rows above them.
$numericRange = "G" + $startRow + ":K" + $lastRow
$ws.Range($numericRange).NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 0; $c -lt $data.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $data[$c]
    }
}
